$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '69.807.91'
Set-TextCell 2 5 '  -1.44%  '
Set-TextCell 3 4 '3.559.89'
Set-TextCell 3 5 '  -2.53%  '
Set-TextCell 4 4 '0.999'
Set-TextCell 4 5 '  -0.09%  '
Set-TextCell 5 4 '574.05'
Set-TextCell 5 5 '  -3.72%  '
Set-TextCell 6 4 '184.98'
Set-TextCell 6 5 '  -4.69%  '
Set-TextCell 7 4 '3.550.91'
Set-TextCell 7 5 '  -2.68%  '
Set-TextCell 8 4 '0.617'
Set-TextCell 8 5 '  -4.69%  '
Set-TextCell 9 5 '  +0.07%  '
Set-TextCell 10 4 '0.182'
Set-TextCell 10 5 '  -0.55%  '
Set-TextCell 11 4 '0.645'
Set-TextCell 11 5 '  -4.21%  '
Set-TextCell 12 4 '54.61'
Set-TextCell 12 5 '  -6.14%  '
Set-TextCell 13 4 '0.0000300'
Set-TextCell 13 5 '  +1.88%  '
Set-TextCell 14 4 '9.48'
Set-TextCell 14 5 '  -4.57%  '
Set-TextCell 15 4 '4.124.81'
Set-TextCell 15 5 '  -2.59%  '
Set-TextCell 16 4 '19.49'
Set-TextCell 16 5 '  -2.85%  '
Set-TextCell 17 4 '3.549.22'
Set-TextCell 17 5 '  -2.72%  '
Set-TextCell 18 4 '69.733.27'
Set-TextCell 18 5 '  -1.50%  '
Set-TextCell 19 4 '12.46'
Set-TextCell 19 5 '  -2.67%  '
Set-TextCell 20 5 '  -0.99%  '
Set-TextCell 21 5 '  -3.93%  '
Set-TextCell 22 4 '490.19'
Set-TextCell 22 5 '  +0.10%  '
Set-TextCell 23 4 '19.08'
Set-TextCell 23 5 '  -0.44%  '
Set-TextCell 24 4 '4.86'
Set-TextCell 24 5 '  -8.13%  '
Set-TextCell 25 4 '4.35'
Set-TextCell 25 5 '  -3.86%  '
Set-TextCell 26 4 '95.01'
Set-TextCell 26 5 '  +3.77%  '
Set-TextCell 27 4 '11.35'
Set-TextCell 27 5 '  -1.08%  '
Set-TextCell 28 4 '2.93'
Set-TextCell 28 5 '  -7.19%  '
Set-TextCell 29 4 '9.23'
Set-TextCell 29 5 '  -3.85%  '
Set-TextCell 30 4 '31.48'
Set-TextCell 30 5 '  -4.23%  '
Set-TextCell 31 4 '7.46'
Set-TextCell 31 5 '  -3.89%  '
Set-TextCell 32 4 '66.37'
Set-TextCell 32 5 '  +0.18%  '
Set-TextCell 33 4 '11.97'
Set-TextCell 33 5 '  -2.61%  '
Set-TextCell 34 4 '0.114'
Set-TextCell 34 5 '  -6.50%  '
Set-TextCell 35 4 '564.87'
Set-TextCell 35 5 '  -10.23%  '
Set-TextCell 36 5 '  +11.78%  '
Set-TextCell 37 4 '38.59'
Set-TextCell 37 5 '  -4.07%  '
Set-TextCell 38 4 '1.00'
Set-TextCell 38 5 '  -0.04%  '
Set-TextCell 39 4 '0.393'
Set-TextCell 39 5 '  -4.90%  '
Set-TextCell 40 4 '0.0₃0785'
Set-TextCell 40 5 '  -4.85%  '
Set-TextCell 41 4 '3.48'
Set-TextCell 41 5 '  -2.81%  '
Set-TextCell 42 4 '3.15'
Set-TextCell 42 5 '  +3.55%  '
Set-TextCell 43 4 '0.133'
Set-TextCell 43 5 '  -10.54%  '
Set-TextCell 44 4 '2.97'
Set-TextCell 44 5 '  -6.06%  '
Set-TextCell 45 4 '3.204.09'
Set-TextCell 45 5 '  -3.07%  '
Set-TextCell 46 5 '  +3.81%  '
Set-TextCell 47 4 '0.0435'
Set-TextCell 47 5 '  -5.02%  '
Set-TextCell 48 4 '9.47'
Set-TextCell 48 5 '  +0.54%  '
Set-TextCell 49 4 '0.135'
Set-TextCell 49 5 '  -2.94%  '
Set-TextCell 50 4 '0.997'
Set-TextCell 50 5 '  -0.18%  '
Set-TextCell 51 4 '3.11'
Set-TextCell 51 5 '  -5.70%  '
